# Update "想去人数" (interest count) figures across the workbook sheets
# to reflect the latest scrape (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 12794
$ws1.Range("F3").Value  = 7156
$ws1.Range("F10").Value = 1003
$ws1.Range("F11").Value = 144
$ws1.Range("F12").Value = 355
$ws1.Range("F13").Value = 1012
$ws1.Range("F18").Value = 245
$ws1.Range("F19").Value = 369
$ws1.Range("F21").Value = 277
$ws1.Range("F22").Value = 310
$ws1.Range("F23").Value = 49
$ws1.Range("F24").Value = 159
$ws1.Range("F25").Value = 371
$ws1.Range("F26").Value = 5233
$ws1.Range("F28").Value = 1427
$ws1.Range("F29").Value = 309
$ws1.Range("F30").Value = 1355
$ws1.Range("F31").Value = 62
$ws1.Range("F32").Value = 37
$ws1.Range("F33").Value = 1361
$ws1.Range("F36").Value = 593
$ws1.Range("F38").Value = 3732

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 3740
$ws2.Range("F5").Value = 3740
$ws2.Range("F8").Value = 54

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9278
$ws3.Range("F3").Value = 560
$ws3.Range("F4").Value = 2011

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9278
$ws4.Range("F3").Value  = 560
$ws4.Range("F4").Value  = 2011
$ws4.Range("F5").Value  = 12794
$ws4.Range("F6").Value  = 7156
$ws4.Range("F8").Value  = 3740
$ws4.Range("F10").Value = 1003
$ws4.Range("F11").Value = 144
$ws4.Range("F12").Value = 355
$ws4.Range("F13").Value = 1012
$ws4.Range("F18").Value = 245
$ws4.Range("F19").Value = 369
$ws4.Range("F21").Value = 277
$ws4.Range("F22").Value = 310
$ws4.Range("F23").Value = 49
$ws4.Range("F27").Value = 159
$ws4.Range("F28").Value = 371
$ws4.Range("F29").Value = 5233
$ws4.Range("F31").Value = 1427
$ws4.Range("F34").Value = 309
$ws4.Range("F36").Value = 1355
$ws4.Range("F37").Value = 62
$ws4.Range("F38").Value = 1361
$ws4.Range("F40").Value = 593
$ws4.Range("F47").Value = 3732
